$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(248,1).Value = 39.78
$ws.Cells.Item(248,2).Value = 0.16
$ws.Cells.Item(248,3).Value = 23.70307803153992
$ws.Cells.Item(248,4).Value = 17.171875

$ws.Cells.Item(249,1).Value = 37
$ws.Cells.Item(249,2).Value = 0.01
$ws.Cells.Item(249,3).Value = 3.178343534469604
$ws.Cells.Item(249,4).Value = 1.640625

$ws.Cells.Item(250,1).Value = 50
$ws.Cells.Item(250,2).Value = 0.0016
$ws.Cells.Item(250,3).Value = 2.165600299835205
$ws.Cells.Item(250,4).Value = 1.40625

$ws.Cells.Item(251,1).Value = 50
$ws.Cells.Item(251,2).Value = 0.0016
$ws.Cells.Item(251,3).Value = 12.8171751499176
$ws.Cells.Item(251,4).Value = 1.453125

$ws.Cells.Item(252,1).Value = 50
$ws.Cells.Item(252,2).Value = 0.0016
$ws.Cells.Item(252,3).Value = 1.738581657409668
$ws.Cells.Item(252,4).Value = 0.53125

$ws.Cells.Item(253,1).Value = 75
$ws.Cells.Item(253,2).Value = 0.0004
$ws.Cells.Item(253,3).Value = 74.63529062271118
$ws.Cells.Item(253,4).Value = 0.828125

$ws.Cells.Item(254,1).Value = 39.44
$ws.Cells.Item(254,2).Value = 0.16
$ws.Cells.Item(254,3).Value = 30.566237449646
$ws.Cells.Item(254,4).Value = 20.5625

$ws.Cells.Item(255,1).Value = 39.44
$ws.Cells.Item(255,2).Value = 0.16
$ws.Cells.Item(255,3).Value = 23.24444532394409
$ws.Cells.Item(255,4).Value = 18.25

$ws.Cells.Item(256,1).Value = 75
$ws.Cells.Item(256,2).Value = 0.0004
$ws.Cells.Item(256,3).Value = 1.669373512268066
$ws.Cells.Item(256,4).Value = 0.53125

$ws.Cells.Item(257,1).Value = 39.44
$ws.Cells.Item(257,2).Value = 0.16
$ws.Cells.Item(257,3).Value = 21.29177212715149
$ws.Cells.Item(257,4).Value = 15.109375

$ws.Cells.Item(258,1).Value = 39.44
$ws.Cells.Item(258,2).Value = 0.16
$ws.Cells.Item(258,3).Value = 21.79749274253845
$ws.Cells.Item(258,4).Value = 14.71875

$ws.Cells.Item(259,1).Value = 50
$ws.Cells.Item(259,2).Value = 0.0016
$ws.Cells.Item(259,3).Value = 2.104458570480347
$ws.Cells.Item(259,4).Value = 0.875

$ws.Cells.Item(260,1).Value = 75
$ws.Cells.Item(260,2).Value = 0.0004
$ws.Cells.Item(260,3).Value = 1.485002517700195
$ws.Cells.Item(260,4).Value = 0.46875

$ws.Cells.Item(261,1).Value = 33.33
$ws.Cells.Item(261,2).Value = 0.0036
$ws.Cells.Item(261,3).Value = 1.890873670578003
$ws.Cells.Item(261,4).Value = 0.9375

$ws.Cells.Item(262,1).Value = 41.5
$ws.Cells.Item(262,2).Value = 0.04000000000000001
$ws.Cells.Item(262,3).Value = 7.52288556098938
$ws.Cells.Item(262,4).Value = 4.3125

$ws.Cells.Item(263,1).Value = 41.5
$ws.Cells.Item(263,2).Value = 0.04000000000000001
$ws.Cells.Item(263,3).Value = 7.478296279907227
$ws.Cells.Item(263,4).Value = 3.484375

$ws.Cells.Item(264,1).Value = 41.5
$ws.Cells.Item(264,2).Value = 0.04000000000000001
$ws.Cells.Item(264,3).Value = 9.29134726524353
$ws.Cells.Item(264,4).Value = 4.53125

$ws.Cells.Item(265,1).Value = 39.44
$ws.Cells.Item(265,2).Value = 0.16
$ws.Cells.Item(265,3).Value = 53.53329634666443
$ws.Cells.Item(265,4).Value = 19.53125

$ws.Cells.Item(266,1).Value = 26.12
$ws.Cells.Item(266,2).Value = 0.16
$ws.Cells.Item(266,3).Value = 48.22412586212158
$ws.Cells.Item(266,4).Value = 17.203125

$ws.Cells.Item(267,1).Value = 40.31
$ws.Cells.Item(267,2).Value = 0.0196
$ws.Cells.Item(267,3).Value = 4.094011783599854
$ws.Cells.Item(267,4).Value = 1.734375

$ws.Cells.Item(268,1).Value = 0
$ws.Cells.Item(268,2).Value = 0.0036
$ws.Cells.Item(268,3).Value = 2.080069303512573
$ws.Cells.Item(268,4).Value = 0.75

$ws.Cells.Item(269,1).Value = 0
$ws.Cells.Item(269,2).Value = 0.0036
$ws.Cells.Item(269,3).Value = 2.129308938980103
$ws.Cells.Item(269,4).Value = 0.734375

$ws.Cells.Item(270,1).Value = 0
$ws.Cells.Item(270,2).Value = 0.04000000000000001
$ws.Cells.Item(270,3).Value = 8.243014097213745
$ws.Cells.Item(270,4).Value = 4.21875

$ws.Cells.Item(271,1).Value = 0
$ws.Cells.Item(271,2).Value = 0.0016
$ws.Cells.Item(271,3).Value = 1.666461706161499
$ws.Cells.Item(271,4).Value = 0.84375

$ws.Cells.Item(272,1).Value = 0
$ws.Cells.Item(272,2).Value = 0.0016
$ws.Cells.Item(272,3).Value = 1.679929733276367
$ws.Cells.Item(272,4).Value = 0.921875

$ws.Cells.Item(273,1).Value = 0
$ws.Cells.Item(273,2).Value = 0.0016
$ws.Cells.Item(273,3).Value = 1.69499397277832
$ws.Cells.Item(273,4).Value = 0.828125

$ws.Cells.Item(274,1).Value = 0
$ws.Cells.Item(274,2).Value = 0.0016
$ws.Cells.Item(274,3).Value = 2.639148950576782
$ws.Cells.Item(274,4).Value = 1.453125

$ws.Cells.Item(275,1).Value = 47
$ws.Cells.Item(275,2).Value = 0.0016
$ws.Cells.Item(275,3).Value = 1.721790790557861
$ws.Cells.Item(275,4).Value = 0.796875

$ws.Cells.Item(276,1).Value = 33.25
$ws.Cells.Item(276,2).Value = 0.16
$ws.Cells.Item(276,3).Value = 23.09088778495789
$ws.Cells.Item(276,4).Value = 17.09375

$ws.Cells.Item(277,1).Value = 33.33
$ws.Cells.Item(277,2).Value = 0.04000000000000001
$ws.Cells.Item(277,3).Value = 7.165266036987305
$ws.Cells.Item(277,4).Value = 4.140625

$ws.Cells.Item(278,1).Value = 82.29000000000001
$ws.Cells.Item(278,2).Value = 0.16
$ws.Cells.Item(278,3).Value = 30.96109223365784
$ws.Cells.Item(278,4).Value = 20.359375

$ws.Cells.Item(279,1).Value = 33.25
$ws.Cells.Item(279,2).Value = 0.16
$ws.Cells.Item(279,3).Value = 29.04900670051575
$ws.Cells.Item(279,4).Value = 18.734375

$ws.Cells.Item(280,1).Value = 50.69
$ws.Cells.Item(280,2).Value = 9
$ws.Cells.Item(280,3).Value = 3794.463066339493
$ws.Cells.Item(280,4).Value = 1564.6875

$ws.Cells.Item(281,1).Value = 0.08
$ws.Cells.Item(281,2).Value = 0.04000000000000001
$ws.Cells.Item(281,3).Value = 7.809202671051025
$ws.Cells.Item(281,4).Value = 4.234375

$ws.Cells.Item(282,1).Value = 0.99
$ws.Cells.Item(282,2).Value = 0.16
$ws.Cells.Item(282,3).Value = 20.09047985076904
$ws.Cells.Item(282,4).Value = 14.078125

$ws.Cells.Item(283,1).Value = 5.13
$ws.Cells.Item(283,2).Value = 0.09
$ws.Cells.Item(283,3).Value = 12.30720639228821
$ws.Cells.Item(283,4).Value = 8.640625

$ws.Cells.Item(284,1).Value = 2.83
$ws.Cells.Item(284,2).Value = 0.16
$ws.Cells.Item(284,3).Value = 20.97641086578369
$ws.Cells.Item(284,4).Value = 16

$ws.Cells.Item(285,1).Value = 16.7
$ws.Cells.Item(285,2).Value = 0.16
$ws.Cells.Item(285,3).Value = 23.43800830841064
$ws.Cells.Item(285,4).Value = 14.078125

$ws.Cells.Item(286,1).Value = 12.21
$ws.Cells.Item(286,2).Value = 0.16
$ws.Cells.Item(286,3).Value = 21.96702885627747
$ws.Cells.Item(286,4).Value = 16.3125

$ws.Cells.Item(287,1).Value = 12.21
$ws.Cells.Item(287,2).Value = 0.16
$ws.Cells.Item(287,3).Value = 26.15172576904297
$ws.Cells.Item(287,4).Value = 15.15625

$ws.Cells.Item(288,1).Value = 12.21
$ws.Cells.Item(288,2).Value = 0.16
$ws.Cells.Item(288,3).Value = 71.19786858558655
$ws.Cells.Item(288,4).Value = 11.65625

$ws.Cells.Item(289,1).Value = 12.21
$ws.Cells.Item(289,2).Value = 0.16
$ws.Cells.Item(289,3).Value = 43.66991925239563
$ws.Cells.Item(289,4).Value = 9.484375

$ws.Cells.Item(290,1).Value = 12.21
$ws.Cells.Item(290,2).Value = 0.16
$ws.Cells.Item(290,3).Value = 48.34774827957153
$ws.Cells.Item(290,4).Value = 11.34375

$ws.Cells.Item(291,1).Value = 12.21
$ws.Cells.Item(291,2).Value = 0.16
$ws.Cells.Item(291,3).Value = 45.7321879863739
$ws.Cells.Item(291,4).Value = 10.15625

$ws.Cells.Item(292,1).Value = 12.21
$ws.Cells.Item(292,2).Value = 0.16
$ws.Cells.Item(292,3).Value = 57.39707326889038
$ws.Cells.Item(292,4).Value = 25.328125

$ws.Cells.Item(293,1).Value = 12.21
$ws.Cells.Item(293,2).Value = 0.16
$ws.Cells.Item(293,3).Value = 66.89149522781372
$ws.Cells.Item(293,4).Value = 28.359375

$ws.Cells.Item(294,1).Value = 12.21
$ws.Cells.Item(294,2).Value = 0.16
$ws.Cells.Item(294,3).Value = 41.22408008575439
$ws.Cells.Item(294,4).Value = 18.3125
